# Update gh-pages to output generated at 456a3b4
# Refresh 展览 (sheet1) and 全部类型 (sheet4) with latest scraped event data.
# 展览: 3 stale/cancelled events removed (IE动漫嘉年华, 书香璃樱动漫游戏嘉年华, 曙光次元动漫游戏嘉年华),
# remaining rows shift up, interest counts ("想去人数") refreshed from source.
# 全部类型: mirrors 展览's changes for the exhibition rows, while the already-unique
# 演出 (performance) rows shift up to follow immediately after.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws4 = $wb.Worksheets.Item(4)

$ws1.Cells.Item(2, 1).Value = 1
$ws1.Cells.Item(2, 2).Value = '2024-09-15'
$ws1.Cells.Item(2, 3).Value = '合肥·ACGN夏日游园会第七回-泳池派对'
$ws1.Cells.Item(2, 4).Value = '金牛路金水里文化产业园 水善汇(金牛路店)'
$ws1.Cells.Item(2, 5).Value = '2024.09.15 09:30-09.16 17:30'
$ws1.Cells.Item(2, 6).Value = 1214
$ws1.Cells.Item(2, 7).Value = 60
$ws1.Cells.Item(2, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91677'
$ws1.Cells.Item(2, 9).Value = '//i2.hdslb.com/bfs/openplatform/202409/j9oW4hzR1725183897413.jpeg'

$ws1.Cells.Item(3, 1).Value = 2
$ws1.Cells.Item(3, 2).Value = '2024-09-15'
$ws1.Cells.Item(3, 3).Value = '蜀山·银泰百货高新店-2024漫趣地带嘉年华（免费）'
$ws1.Cells.Item(3, 4).Value = '高新区望江西路888号 银泰百货（高新店）'
$ws1.Cells.Item(3, 5).Value = '2024.09.15 10:00-10.02 22:00'
$ws1.Cells.Item(3, 6).Value = 227
$ws1.Cells.Item(3, 7).Value = 30
$ws1.Cells.Item(3, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91869'
$ws1.Cells.Item(3, 9).Value = '//i2.hdslb.com/bfs/openplatform/202409/JDGIWMyd1725422862878.png'

$ws1.Cells.Item(4, 1).Value = 3
$ws1.Cells.Item(4, 2).Value = '2024-09-16'
$ws1.Cells.Item(4, 3).Value = '肥西·星域动漫游戏嘉年华'
$ws1.Cells.Item(4, 4).Value = '金寨路与云谷路交口金云国际9号楼商(邮政银行旁边) 吉祥如意大酒店(肥西店)'
$ws1.Cells.Item(4, 5).Value = '2024.09.16 10:00-09.16 17:00'
$ws1.Cells.Item(4, 6).Value = 67
$ws1.Cells.Item(4, 7).Value = 45
$ws1.Cells.Item(4, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90489'
$ws1.Cells.Item(4, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/6xk6G8E71722525186252.jpeg'

$ws1.Cells.Item(5, 1).Value = 4
$ws1.Cells.Item(5, 2).Value = '2024-09-21'
$ws1.Cells.Item(5, 3).Value = '合肥·漫有引力动漫游戏嘉年华'
$ws1.Cells.Item(5, 4).Value = '幸福路1号(筑梦集团·结婚产业园·B1幢) 费加罗宴会艺术中心(旗舰店)'
$ws1.Cells.Item(5, 5).Value = '2024.09.21 10:00-09.21 17:00'
$ws1.Cells.Item(5, 6).Value = 67
$ws1.Cells.Item(5, 7).Value = 50
$ws1.Cells.Item(5, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90905'
$ws1.Cells.Item(5, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/gfeOndjM1723659151069.png'

$ws1.Cells.Item(6, 1).Value = 5
$ws1.Cells.Item(6, 2).Value = '2024-10-01'
$ws1.Cells.Item(6, 3).Value = '合肥·星域动漫游戏嘉年华'
$ws1.Cells.Item(6, 4).Value = '新站区东方大道288号 少荃体育中心'
$ws1.Cells.Item(6, 5).Value = '2024.10.01 10:00-10.01 17:00'
$ws1.Cells.Item(6, 6).Value = 11
$ws1.Cells.Item(6, 7).Value = 58
$ws1.Cells.Item(6, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91878'
$ws1.Cells.Item(6, 9).Value = '//i0.hdslb.com/bfs/openplatform/202409/NOg6Wwjh1725121441581.png'

$ws1.Cells.Item(7, 1).Value = 6
$ws1.Cells.Item(7, 2).Value = '2024-10-01'
$ws1.Cells.Item(7, 3).Value = '合肥·第十五届次元之门动漫游戏博览会'
$ws1.Cells.Item(7, 4).Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$ws1.Cells.Item(7, 5).Value = '2024.10.01 09:30-10.02 17:30'
$ws1.Cells.Item(7, 6).Value = 5698
$ws1.Cells.Item(7, 7).Value = 70
$ws1.Cells.Item(7, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91133'
$ws1.Cells.Item(7, 9).Value = '//i1.hdslb.com/bfs/openplatform/202408/PlcqtYWR1724315434068.jpeg'

$ws1.Cells.Item(8, 1).Value = 7
$ws1.Cells.Item(8, 2).Value = '2024-10-01'
$ws1.Cells.Item(8, 3).Value = '合肥·首届AT次元时代动漫游戏嘉年华'
$ws1.Cells.Item(8, 4).Value = '凤淮路与公园路交叉口南行50米路西 庐阳区全民健身中心'
$ws1.Cells.Item(8, 5).Value = '2024.10.01 09:30-10.03 17:00'
$ws1.Cells.Item(8, 6).Value = 5049
$ws1.Cells.Item(8, 7).Value = 68
$ws1.Cells.Item(8, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90908'
$ws1.Cells.Item(8, 9).Value = '//i0.hdslb.com/bfs/openplatform/202408/Y6P7qrm11724139193256.jpeg'

$ws1.Cells.Item(9, 1).Value = 8
$ws1.Cells.Item(9, 2).Value = '2024-10-04'
$ws1.Cells.Item(9, 3).Value = '合肥·Holic动漫游戏展'
$ws1.Cells.Item(9, 4).Value = '庐州大道800号 合肥融创茂'
$ws1.Cells.Item(9, 5).Value = '2024.10.04 10:00-10.06 17:00'
$ws1.Cells.Item(9, 6).Value = 23
$ws1.Cells.Item(9, 7).Value = 55
$ws1.Cells.Item(9, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92061'
$ws1.Cells.Item(9, 9).Value = '//i1.hdslb.com/bfs/openplatform/202409/AZ0LsUce1725522015668.jpeg'

$ws1.Cells.Item(10, 1).Value = 9
$ws1.Cells.Item(10, 2).Value = '2024-10-04'
$ws1.Cells.Item(10, 3).Value = '合肥·乐帮•崩原铁绝only同人首展'
$ws1.Cells.Item(10, 4).Value = '丹霞路488号金星商业城三楼 迷鹿轰趴'
$ws1.Cells.Item(10, 5).Value = '2024.10.04 10:00-10.05 16:30'
$ws1.Cells.Item(10, 6).Value = 52
$ws1.Cells.Item(10, 7).Value = 58
$ws1.Cells.Item(10, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91524'
$ws1.Cells.Item(10, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/739I7YRr1724912450704.png'

$ws1.Cells.Item(11, 1).Value = 10
$ws1.Cells.Item(11, 2).Value = '2024-10-06'
$ws1.Cells.Item(11, 3).Value = '合肥·星月动漫游戏展'
$ws1.Cells.Item(11, 4).Value = '灵石路与皇藏峪路交叉口西南10米安徽百事兴电气有限公司院内2栋厂房2层 兄弟篮球俱乐部'
$ws1.Cells.Item(11, 5).Value = '2024.10.06 10:00-10.06 17:00'
$ws1.Cells.Item(11, 6).Value = 5
$ws1.Cells.Item(11, 7).Value = 45
$ws1.Cells.Item(11, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91958'
$ws1.Cells.Item(11, 9).Value = '//i2.hdslb.com/bfs/openplatform/202409/mgB8U6bN1725361649767.jpeg'

$ws1.Cells.Item(12, 1).Value = 11
$ws1.Cells.Item(12, 2).Value = '2024-10-06'
$ws1.Cells.Item(12, 3).Value = '合肥·首届火影忍者同人only'
$ws1.Cells.Item(12, 4).Value = '长江东路金太阳家具广场南门二楼 优极篮球馆'
$ws1.Cells.Item(12, 5).Value = '2024.10.06 09:30-10.06 17:30'
$ws1.Cells.Item(12, 6).Value = 58
$ws1.Cells.Item(12, 7).Value = 75
$ws1.Cells.Item(12, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91658'
$ws1.Cells.Item(12, 9).Value = '//i0.hdslb.com/bfs/openplatform/202408/f8ylbskH1725027552569.jpeg'

$ws1.Cells.Item(13, 1).Value = 12
$ws1.Cells.Item(13, 2).Value = '2024-10-26'
$ws1.Cells.Item(13, 3).Value = '合肥·W·A第五人格同人only2.0'
$ws1.Cells.Item(13, 4).Value = '莲花路与石门路交口西北角（尚泽大都会B座四楼） 格律诗婚礼艺术中心(经开店)'
$ws1.Cells.Item(13, 5).Value = '2024.10.26 09:30-10.26 17:00'
$ws1.Cells.Item(13, 6).Value = 208
$ws1.Cells.Item(13, 7).Value = 68
$ws1.Cells.Item(13, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91123'
$ws1.Cells.Item(13, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/YqXHTFM81724066565119.png'

$ws1.Cells.Item(14, 1).Value = 13
$ws1.Cells.Item(14, 2).Value = '2024-11-17'
$ws1.Cells.Item(14, 3).Value = '合肥·MAX特摄同人only2.0'
$ws1.Cells.Item(14, 4).Value = '桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间'
$ws1.Cells.Item(14, 5).Value = '2024.11.17 10:00-11.17 18:00'
$ws1.Cells.Item(14, 6).Value = 11
$ws1.Cells.Item(14, 7).Value = 60
$ws1.Cells.Item(14, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92498'
$ws1.Cells.Item(14, 9).Value = '//i0.hdslb.com/bfs/openplatform/202409/R4WJxlQe1726230330813.jpeg'

# Remove now-stale trailing rows 15-17 (old data beyond the new 13-row table)
$ws1.Range("A15:I17").Delete()

$ws4.Cells.Item(2, 1).Value = 1
$ws4.Cells.Item(2, 2).Value = '2024-09-15'
$ws4.Cells.Item(2, 3).Value = '合肥·ACGN夏日游园会第七回-泳池派对'
$ws4.Cells.Item(2, 4).Value = '金牛路金水里文化产业园 水善汇(金牛路店)'
$ws4.Cells.Item(2, 5).Value = '2024.09.15 09:30-09.16 17:30'
$ws4.Cells.Item(2, 6).Value = 1214
$ws4.Cells.Item(2, 7).Value = 60
$ws4.Cells.Item(2, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91677'
$ws4.Cells.Item(2, 9).Value = '//i2.hdslb.com/bfs/openplatform/202409/j9oW4hzR1725183897413.jpeg'

$ws4.Cells.Item(3, 1).Value = 2
$ws4.Cells.Item(3, 2).Value = '2024-09-15'
$ws4.Cells.Item(3, 3).Value = '蜀山·银泰百货高新店-2024漫趣地带嘉年华（免费）'
$ws4.Cells.Item(3, 4).Value = '高新区望江西路888号 银泰百货（高新店）'
$ws4.Cells.Item(3, 5).Value = '2024.09.15 10:00-10.02 22:00'
$ws4.Cells.Item(3, 6).Value = 227
$ws4.Cells.Item(3, 7).Value = 30
$ws4.Cells.Item(3, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91869'
$ws4.Cells.Item(3, 9).Value = '//i2.hdslb.com/bfs/openplatform/202409/JDGIWMyd1725422862878.png'

$ws4.Cells.Item(4, 1).Value = 3
$ws4.Cells.Item(4, 2).Value = '2024-09-16'
$ws4.Cells.Item(4, 3).Value = '肥西·星域动漫游戏嘉年华'
$ws4.Cells.Item(4, 4).Value = '金寨路与云谷路交口金云国际9号楼商(邮政银行旁边) 吉祥如意大酒店(肥西店)'
$ws4.Cells.Item(4, 5).Value = '2024.09.16 10:00-09.16 17:00'
$ws4.Cells.Item(4, 6).Value = 67
$ws4.Cells.Item(4, 7).Value = 45
$ws4.Cells.Item(4, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90489'
$ws4.Cells.Item(4, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/6xk6G8E71722525186252.jpeg'

$ws4.Cells.Item(5, 1).Value = 4
$ws4.Cells.Item(5, 2).Value = '2024-09-21'
$ws4.Cells.Item(5, 3).Value = '合肥·漫有引力动漫游戏嘉年华'
$ws4.Cells.Item(5, 4).Value = '幸福路1号(筑梦集团·结婚产业园·B1幢) 费加罗宴会艺术中心(旗舰店)'
$ws4.Cells.Item(5, 5).Value = '2024.09.21 10:00-09.21 17:00'
$ws4.Cells.Item(5, 6).Value = 67
$ws4.Cells.Item(5, 7).Value = 50
$ws4.Cells.Item(5, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90905'
$ws4.Cells.Item(5, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/gfeOndjM1723659151069.png'

$ws4.Cells.Item(6, 1).Value = 5
$ws4.Cells.Item(6, 2).Value = '2024-10-01'
$ws4.Cells.Item(6, 3).Value = '合肥·星域动漫游戏嘉年华'
$ws4.Cells.Item(6, 4).Value = '新站区东方大道288号 少荃体育中心'
$ws4.Cells.Item(6, 5).Value = '2024.10.01 10:00-10.01 17:00'
$ws4.Cells.Item(6, 6).Value = 11
$ws4.Cells.Item(6, 7).Value = 58
$ws4.Cells.Item(6, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91878'
$ws4.Cells.Item(6, 9).Value = '//i0.hdslb.com/bfs/openplatform/202409/NOg6Wwjh1725121441581.png'

$ws4.Cells.Item(7, 1).Value = 6
$ws4.Cells.Item(7, 2).Value = '2024-10-01'
$ws4.Cells.Item(7, 3).Value = '合肥·第十五届次元之门动漫游戏博览会'
$ws4.Cells.Item(7, 4).Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$ws4.Cells.Item(7, 5).Value = '2024.10.01 09:30-10.02 17:30'
$ws4.Cells.Item(7, 6).Value = 5698
$ws4.Cells.Item(7, 7).Value = 70
$ws4.Cells.Item(7, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91133'
$ws4.Cells.Item(7, 9).Value = '//i1.hdslb.com/bfs/openplatform/202408/PlcqtYWR1724315434068.jpeg'

$ws4.Cells.Item(8, 1).Value = 7
$ws4.Cells.Item(8, 2).Value = '2024-10-01'
$ws4.Cells.Item(8, 3).Value = '合肥·首届AT次元时代动漫游戏嘉年华'
$ws4.Cells.Item(8, 4).Value = '凤淮路与公园路交叉口南行50米路西 庐阳区全民健身中心'
$ws4.Cells.Item(8, 5).Value = '2024.10.01 09:30-10.03 17:00'
$ws4.Cells.Item(8, 6).Value = 5049
$ws4.Cells.Item(8, 7).Value = 68
$ws4.Cells.Item(8, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90908'
$ws4.Cells.Item(8, 9).Value = '//i0.hdslb.com/bfs/openplatform/202408/Y6P7qrm11724139193256.jpeg'

$ws4.Cells.Item(9, 1).Value = 8
$ws4.Cells.Item(9, 2).Value = '2024-10-04'
$ws4.Cells.Item(9, 3).Value = '合肥·Holic动漫游戏展'
$ws4.Cells.Item(9, 4).Value = '庐州大道800号 合肥融创茂'
$ws4.Cells.Item(9, 5).Value = '2024.10.04 10:00-10.06 17:00'
$ws4.Cells.Item(9, 6).Value = 23
$ws4.Cells.Item(9, 7).Value = 55
$ws4.Cells.Item(9, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92061'
$ws4.Cells.Item(9, 9).Value = '//i1.hdslb.com/bfs/openplatform/202409/AZ0LsUce1725522015668.jpeg'

$ws4.Cells.Item(10, 1).Value = 9
$ws4.Cells.Item(10, 2).Value = '2024-10-04'
$ws4.Cells.Item(10, 3).Value = '合肥·乐帮•崩原铁绝only同人首展'
$ws4.Cells.Item(10, 4).Value = '丹霞路488号金星商业城三楼 迷鹿轰趴'
$ws4.Cells.Item(10, 5).Value = '2024.10.04 10:00-10.05 16:30'
$ws4.Cells.Item(10, 6).Value = 52
$ws4.Cells.Item(10, 7).Value = 58
$ws4.Cells.Item(10, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91524'
$ws4.Cells.Item(10, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/739I7YRr1724912450704.png'

$ws4.Cells.Item(11, 1).Value = 10
$ws4.Cells.Item(11, 2).Value = '2024-10-06'
$ws4.Cells.Item(11, 3).Value = '合肥·星月动漫游戏展'
$ws4.Cells.Item(11, 4).Value = '灵石路与皇藏峪路交叉口西南10米安徽百事兴电气有限公司院内2栋厂房2层 兄弟篮球俱乐部'
$ws4.Cells.Item(11, 5).Value = '2024.10.06 10:00-10.06 17:00'
$ws4.Cells.Item(11, 6).Value = 5
$ws4.Cells.Item(11, 7).Value = 45
$ws4.Cells.Item(11, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91958'
$ws4.Cells.Item(11, 9).Value = '//i2.hdslb.com/bfs/openplatform/202409/mgB8U6bN1725361649767.jpeg'

$ws4.Cells.Item(12, 1).Value = 11
$ws4.Cells.Item(12, 2).Value = '2024-10-06'
$ws4.Cells.Item(12, 3).Value = '合肥·首届火影忍者同人only'
$ws4.Cells.Item(12, 4).Value = '长江东路金太阳家具广场南门二楼 优极篮球馆'
$ws4.Cells.Item(12, 5).Value = '2024.10.06 09:30-10.06 17:30'
$ws4.Cells.Item(12, 6).Value = 58
$ws4.Cells.Item(12, 7).Value = 75
$ws4.Cells.Item(12, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91658'
$ws4.Cells.Item(12, 9).Value = '//i0.hdslb.com/bfs/openplatform/202408/f8ylbskH1725027552569.jpeg'

$ws4.Cells.Item(13, 1).Value = 12
$ws4.Cells.Item(13, 2).Value = '2024-10-26'
$ws4.Cells.Item(13, 3).Value = '合肥·W·A第五人格同人only2.0'
$ws4.Cells.Item(13, 4).Value = '莲花路与石门路交口西北角（尚泽大都会B座四楼） 格律诗婚礼艺术中心(经开店)'
$ws4.Cells.Item(13, 5).Value = '2024.10.26 09:30-10.26 17:00'
$ws4.Cells.Item(13, 6).Value = 208
$ws4.Cells.Item(13, 7).Value = 68
$ws4.Cells.Item(13, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91123'
$ws4.Cells.Item(13, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/YqXHTFM81724066565119.png'

$ws4.Cells.Item(14, 1).Value = 13
$ws4.Cells.Item(14, 2).Value = '2024-10-26'
$ws4.Cells.Item(14, 3).Value = '合肥·《四月是你的谎言》—“公生”与“薰”的钢琴小提琴唯美经典音乐集'
$ws4.Cells.Item(14, 4).Value = '徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院'
$ws4.Cells.Item(14, 5).Value = '2024.10.26 19:30-10.26 21:00'
$ws4.Cells.Item(14, 6).Value = 78
$ws4.Cells.Item(14, 7).Value = 80
$ws4.Cells.Item(14, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90322'
$ws4.Cells.Item(14, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/BiVgXUKH1722824304648.jpeg'

$ws4.Cells.Item(15, 1).Value = 14
$ws4.Cells.Item(15, 2).Value = '2024-11-09'
$ws4.Cells.Item(15, 3).Value = '合肥·一生必听的钢琴曲—“从巴赫 · 莫扎特到肖邦 · 李斯特”钢琴圣手谭小棠独奏音乐会'
$ws4.Cells.Item(15, 4).Value = '徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院'
$ws4.Cells.Item(15, 5).Value = '2024.11.09 19:30-11.09 21:00'
$ws4.Cells.Item(15, 6).Value = 5
$ws4.Cells.Item(15, 7).Value = 64
$ws4.Cells.Item(15, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90593'
$ws4.Cells.Item(15, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/SYfLxnO21723442234232.jpeg'

$ws4.Cells.Item(16, 1).Value = 15
$ws4.Cells.Item(16, 2).Value = '2024-11-17'
$ws4.Cells.Item(16, 3).Value = '合肥·MAX特摄同人only2.0'
$ws4.Cells.Item(16, 4).Value = '桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间'
$ws4.Cells.Item(16, 5).Value = '2024.11.17 10:00-11.17 18:00'
$ws4.Cells.Item(16, 6).Value = 11
$ws4.Cells.Item(16, 7).Value = 60
$ws4.Cells.Item(16, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92498'
$ws4.Cells.Item(16, 9).Value = '//i0.hdslb.com/bfs/openplatform/202409/R4WJxlQe1726230330813.jpeg'

$ws4.Cells.Item(17, 1).Value = 16
$ws4.Cells.Item(17, 2).Value = '2024-12-07'
$ws4.Cells.Item(17, 3).Value = '合肥·一生必听的古典系列《钟》—超技钢琴曲炫彩音乐会'
$ws4.Cells.Item(17, 4).Value = '徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院'
$ws4.Cells.Item(17, 5).Value = '2024.12.07 19:30-12.07 21:00'
$ws4.Cells.Item(17, 6).Value = 2
$ws4.Cells.Item(17, 7).Value = 64
$ws4.Cells.Item(17, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91608'
$ws4.Cells.Item(17, 9).Value = '//i0.hdslb.com/bfs/openplatform/202408/wiLiWoeM1725005636569.jpeg'

# Remove now-stale trailing rows 18-20 (old data beyond the new 16-row table)
$ws4.Range("A18:I20").Delete()
